# Insert a new data row right before the current row 78 ("Hortaliza, Vega
# Monumental Concepción - Espinaca" sheet). Every existing row from 78
# downward shifts down by one (78->79, ..., 134->135), and the newly
# opened row 78 is populated with a new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 78:134 down to 79:135, opening up a blank row 78.
$ws.Rows("78:78").Insert()

# Populate the new row 78 with the new observation.
$ws.Range("A78").Value = 11
$ws.Range("B78").Value = "Vega Monumental Concepción"
$ws.Range("C78").Value = "Bíobío"
$ws.Range("D78").Value = 45161
$ws.Range("E78").Value = 8
$ws.Range("F78").Value = 100112012
$ws.Range("G78").Value = "Espinaca"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 50
$ws.Range("K78").Value = 6500
$ws.Range("L78").Value = 7000
$ws.Range("M78").Value = 6700
$ws.Range("N78").Value = "$/cuna 10 kilos"
$ws.Range("O78").Value = "Región Metropolitana"
$ws.Range("P78").Value = 670
$ws.Range("Q78").Value = 10
$ws.Range("R78").Value = "Hortaliza"
